$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1. Update the cached "datetimeFigureOut" field text on the slide
#    master and every slide layout from 12/01/2021 to 17/01/2021.
# ---------------------------------------------------------------
$newDate = "17/01/2021"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------
# 2. Resize / reposition "Group 18" on slide 2 (and its two child
#    textboxes "TextBox 16" / "TextBox 17") per the updated layout.
# ---------------------------------------------------------------
$s = $p.Slides.Item(2)
$grp = $s.Shapes.Item(4)

$ptPerEmu = 1.0 / 914400.0 * 72.0

$grp.Left   = -1 * 0
$grp.Top    = -6707 * $ptPerEmu
$grp.Width  = 8101609 * $ptPerEmu
$grp.Height = 3549511 * $ptPerEmu

for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $child = $grp.GroupItems.Item($i)
    if ($child.Name -eq "TextBox 16") {
        $child.Left = 0 * $ptPerEmu
        $child.Top  = -6706 * $ptPerEmu
    } elseif ($child.Name -eq "TextBox 17") {
        $child.Left = 3978085 * $ptPerEmu
        $child.Top  = -6707 * $ptPerEmu
    }
}
